# "edit reference notation in CBM thickness"
# Rewrite the short author/year citation strings in column A of the
# "CBM thickness" sheet into the fuller "<Author> et al., <year>" notation
# (with a few disambiguating qualifiers), and widen column A to fit the
# new, longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")

$ws.Range("A2").Value = "Rodrigues et al., 1983"
$ws.Range("A3").Value = "Ceafalan et al., 2019"
$ws.Range("A4").Value = "Creutzfeldt et al., 1970"
$ws.Range("A6").Value = "Calson et al., 2003 (muscle)"
$ws.Range("A5").Value = "Calson et al., 2003 (retina)"
$ws.Range("A7").Value = "Lash et al., 1989 (11 wk.)"
$ws.Range("A8").Value = "Lash et al., 1989 (18 wk.)"
$ws.Range("A9").Value = "Danis & Yang, 1993"

# Column A needs to be wider now that the labels are longer.
$ws.Columns.Item(1).ColumnWidth = 24.1640625

# Saved selection state.
$ws.Range("B9").Select()
